# Update "想去人数" (interested-people count) values in the "展览" and
# "全部类型" worksheets to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F (想去人数), per sheet.
$updates1 = @{
    2  = 833
    3  = 8
    4  = 1149
    6  = 12396
    9  = 497
    10 = 446
    12 = 916
    13 = 13622
    14 = 13840
    19 = 1038
    22 = 266
    23 = 4935
    24 = 224
}

$updates4 = @{
    2  = 833
    3  = 8
    4  = 1149
    6  = 12396
    9  = 497
    10 = 446
    12 = 916
    13 = 13622
    14 = 13840
    19 = 1038
    22 = 267
    23 = 4935
    24 = 224
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
